$d = $word.ActiveDocument

# The target text currently reads (paraphrased):
#   ... the below msi file located under "angularjs-app\formApp\angular\http-server"
#   and install using all the default options.
# It must become:
#   ... the below msi file located under "http-server" and install using all the
#   default options.
# (path for the httpserver msi simplified) while keeping the "_GoBack" bookmark
# that originally sat between "formApp\" and "angular\" in place (it should end up
# immediately after the opening quote, right before "http-server").

# Step 1: replace everything from the opening quote through the trailing
# "angular\" (just before "http-server") with a placeholder marker character.
# Doing the replace in one shot (rather than leaving the boundary exactly on
# "angular") avoids stray/orphaned w:proofErr spell-check markers being left
# behind in the document body.
$find1 = $d.Content.Find
$found1 = $find1.Execute("`"angularjs-app\formApp\angular\", $true, $false, $false, $false, $false, $true, 1, $false, "@@MARK@@", 2)

# Step 2: locate the placeholder and swap it back for a plain straight quote
# using a direct Range.Text assignment (NOT Find/Replace) so Word's
# smart-quote AutoCorrect never kicks in and mangles it into a curly quote.
$r2 = $d.Content
$f2 = $r2.Find
$found2 = $f2.Execute("@@MARK@@http-server")
$markStart = $r2.Start
$markEnd = $markStart + 8
$markRng = $d.Range($markStart, $markEnd)
$markRng.Text = """"

# Step 3: the run now reads " file located under """ as a single merged run
# (Word merges identically-formatted adjacent runs when text between them is
# deleted). Re-split it back into the original two runs - " file located "
# and "under """ - by nudging the formatting on the "under """ portion (set
# then clear Bold, a value it already has, which forces the run boundary
# without altering the visible formatting).
$r3 = $d.Content
$f3 = $r3.Find
$found3 = $f3.Execute("under `"")
$underStart = $r3.Start
$underEnd = $r3.End
$splitRng = $d.Range($underStart, $underEnd)
$splitRng.Bold = 1
$splitRng.Bold = 0

# Step 4: re-create the "_GoBack" bookmark immediately after the quote
# (right before "http-server"), matching its original relative position.
$bmRng = $d.Range($underEnd, $underEnd)
$d.Bookmarks.Add("_GoBack", $bmRng)
